$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.288.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.521.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "196.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "583.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.33%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.626"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000287"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "679.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +13.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.079.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.358.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.540.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.15%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.122"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.965"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "106.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.33%  "
$ws.Range("E34").Value = "  -4.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "62.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.809.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0812"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "502.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.372"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.74%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.134"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "34.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.53%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +20.00%  "
